$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# DISTRICT column (B) value updates
$ws.Cells.Item(2, 2).Value = 11
$ws.Cells.Item(3, 2).Value = 20
$ws.Cells.Item(4, 2).Value = 110
$ws.Cells.Item(5, 2).Value = 10011

# New Inspection_Date column (C) values - Excel date serials for
# 2022-01-01, 2022-02-01, 2022-03-01, 2022-04-01
$ws.Cells.Item(2, 3).Value = 44562
$ws.Cells.Item(3, 3).Value = 44593
$ws.Cells.Item(4, 3).Value = 44621
$ws.Cells.Item(5, 3).Value = 44652

# Format the first date cell as a short date, then propagate the same
# formatting (and only the formatting) to the rest of the new date cells
# so they all resolve to a single shared style.
$ws.Cells.Item(2, 3).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 3).Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to B6, matching the saved cursor position
$ws.Range("B6").Select()
